$wb = $excel.ActiveWorkbook
$wsRuns = $wb.Worksheets.Item("Runs")

$wsRuns.Range("A1").Value = 0
$wsRuns.Range("B1").Value = "Run105"
$wsRuns.Range("C1").Value = "0.560 (0.449)"
$wsRuns.Range("D1").Value = "0.455 (0.287)"
$wsRuns.Range("E1").Value = "0.600 (0.491)"
$wsRuns.Range("F1").Value = "0.560 (0.428)"
$wsRuns.Range("G1").Value = "0.580 (0.244)"
$wsRuns.Range("H1").Value = "0.550 (0.499)"
$wsRuns.Range("I1").Value = "0.768 (0.386)"
$wsRuns.Range("J1").Value = "0.716 (0.210)"
$wsRuns.Range("K1").Value = "0.779 (0.415)"
$wsRuns.Range("L1").Value = "0.636 (0.435)"
$wsRuns.Range("M1").Value = "0.555 (0.280)"
$wsRuns.Range("N1").Value = "0.664 (0.473)"

$wsRuns.Range("A2").Value = 0
$wsRuns.Range("B2").Value = "Run108"
$wsRuns.Range("C2").Value = "0.509 (0.444)"
$wsRuns.Range("D2").Value = "0.758 (0.160)"
$wsRuns.Range("E2").Value = "0.382 (0.486)"
$wsRuns.Range("F2").Value = "0.617 (0.456)"
$wsRuns.Range("G2").Value = "0.793 (0.197)"
$wsRuns.Range("H2").Value = "0.564 (0.497)"
$wsRuns.Range("I2").Value = "0.759 (0.414)"
$wsRuns.Range("J2").Value = "0.607 (0.324)"
$wsRuns.Range("K2").Value = "0.774 (0.419)"
$wsRuns.Range("L2").Value = "0.613 (0.450)"
$wsRuns.Range("M2").Value = "0.745 (0.202)"
$wsRuns.Range("N2").Value = "0.572 (0.495)"

$wsRuns.Range("A3").Value = 0
$wsRuns.Range("B3").Value = "Run109"
$wsRuns.Range("C3").Value = "0.743 (0.383)"
$wsRuns.Range("D3").Value = "0.733 (0.173)"
$wsRuns.Range("E3").Value = "0.746 (0.436)"
$wsRuns.Range("F3").Value = "0.776 (0.390)"
$wsRuns.Range("G3").Value = "0.885 (0.100)"
$wsRuns.Range("H3").Value = "0.741 (0.440)"
$wsRuns.Range("I3").Value = "0.901 (0.260)"
$wsRuns.Range("J3").Value = "0.801 (0.087)"
$wsRuns.Range("K3").Value = "0.918 (0.275)"
$wsRuns.Range("L3").Value = "0.807 (0.352)"
$wsRuns.Range("M3").Value = "0.779 (0.156)"
$wsRuns.Range("N3").Value = "0.814 (0.389)"

$wsRuns.Range("A4").Value = 0
$wsRuns.Range("B4").Value = "Run110"
$wsRuns.Range("C4").Value = "0.720 (0.439)"
$wsRuns.Range("D4").Value = "0.722 (0.169)"
$wsRuns.Range("E4").Value = "0.719 (0.450)"
$wsRuns.Range("F4").Value = "0.864 (0.299)"
$wsRuns.Range("G4").Value = "0.875 (0.132)"
$wsRuns.Range("H4").Value = "0.859 (0.349)"
$wsRuns.Range("I4").Value = "0.816 (0.341)"
$wsRuns.Range("J4").Value = "0.803 (0.141)"
$wsRuns.Range("K4").Value = "0.820 (0.384)"
$wsRuns.Range("L4").Value = "0.781 (0.386)"
$wsRuns.Range("M4").Value = "0.809 (0.150)"
$wsRuns.Range("N4").Value = "0.775 (0.418)"

$wsRuns.Range("A5").Value = 0
$wsRuns.Range("B5").Value = "Run111"
$wsRuns.Range("C5").Value = "0.775 (0.394)"
$wsRuns.Range("D5").Value = "0.872 (0.074)"
$wsRuns.Range("E5").Value = "0.752 (0.432)"
$wsRuns.Range("F5").Value = "0.810 (0.341)"
$wsRuns.Range("G5").Value = "0.826 (0.159)"
$wsRuns.Range("H5").Value = "0.802 (0.399)"
$wsRuns.Range("I5").Value = "0.879 (0.307)"
$wsRuns.Range("J5").Value = "0.676 (0.220)"
$wsRuns.Range("K5").Value = "0.895 (0.307)"
$wsRuns.Range("L5").Value = "0.819 (0.358)"
$wsRuns.Range("M5").Value = "0.826 (0.152)"
$wsRuns.Range("N5").Value = "0.818 (0.386)"

$wsRuns.Range("A6").Value = 0
$wsRuns.Range("B6").Value = "Run116"
$wsRuns.Range("C6").Value = "0.130 (0.214)"
$wsRuns.Range("D6").Value = "0.399 (0.183)"
$wsRuns.Range("E6").Value = "0.000 (0.000)"
$wsRuns.Range("F6").Value = "0.250 (0.323)"
$wsRuns.Range("G6").Value = "0.531 (0.216)"
$wsRuns.Range("H6").Value = "0.052 (0.223)"
$wsRuns.Range("I6").Value = "0.142 (0.288)"
$wsRuns.Range("J6").Value = "0.549 (0.176)"
$wsRuns.Range("K6").Value = "0.053 (0.224)"
$wsRuns.Range("L6").Value = "0.151 (0.266)"
$wsRuns.Range("M6").Value = "0.466 (0.201)"
$wsRuns.Range("N6").Value = "0.030 (0.172)"

$wsIter = $wb.Worksheets.Item("Iterations")
$wsIter.Range("A1").Value = 0
$wsIter.Range("B1").Value = "Iteration30"
$wsIter.Range("C1").Value = "0.69 (0.10)"
$wsIter.Range("D1").Value = "0.77 (0.06)"
$wsIter.Range("E1").Value = "0.65 (0.16)"
$wsIter.Range("F1").Value = "0.77 (0.09)"
$wsIter.Range("G1").Value = "0.84 (0.04)"
$wsIter.Range("H1").Value = "0.74 (0.11)"
$wsIter.Range("I1").Value = "0.84 (0.06)"
$wsIter.Range("J1").Value = "0.72 (0.08)"
$wsIter.Range("K1").Value = "0.85 (0.06)"
$wsIter.Range("L1").Value = "0.75 (0.08)"
$wsIter.Range("M1").Value = "0.79 (0.03)"
$wsIter.Range("N1").Value = "0.74 (0.10)"
